$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column into the worksheet at column D, shifting existing
# "is null?" / "is key?" / "initial value" / "description" columns right.
$ws.Range("D:D").Insert()

# New header cell for the inserted "type in JavaDB" column.
$ws.Range("D4").Value = "type in JavaDB"

# Fill in the data cells for the new column.
$ws.Range("D5").Value = "INT"
$ws.Range("D6").Value = "VARCHAR(50)"
$ws.Range("D7").Value = "varchar(200)"

# Resize / refresh the Excel table (ListObject) so it covers the new column.
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("B4:H7"))

$ws.Range("D4").Select()
